# Reporting year sheet for 2000: add the year value in the row below the
# header, move the selection onto the new cell, and set the page to
# portrait orientation (matches the saved page-setup state in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2000
$ws.Range("A2").Select()

# xlPortrait = 1
$ws.PageSetup.Orientation = 1
